# Apply the edits described by the diff:
#  - Scroll the sheet view back to the top (remove topLeftCell="A98")
#    and move/collapse the selection to B1 (was C121).
#  - Resize column B from its previous auto "best fit" width to a
#    manually-set custom width of 130 pixels (≈17.857143 characters,
#    which XLSX serializes as width="18.5703125").
#  - Touch the autoFilter range so Excel re-confirms/normalizes it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("city_data")

# Make sure this is the active sheet/window context.
$ws.Activate()

# --- Column B width change -------------------------------------------------
# Previously bestFit (auto-sized) at ~12.29 chars; now a manual custom width.
$ws.Columns.Item(2).ColumnWidth = 125 / 7

# --- Selection / scroll position --------------------------------------------
# Move the active selection to B1 and scroll the view back so A1 is visible
# again (clears the old topLeftCell="A98" / selection at C121).
$ws.Range("A1").Select() | Out-Null
$ws.Range("B1").Select() | Out-Null

# --- Re-apply the AutoFilter over the same range ----------------------------
$ws.Range("A1:C121").AutoFilter(1) | Out-Null

$wb.Save()
